$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.172.97'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -1.75%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.849.01'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -2.59%  '
$ws.Range('E3').ClearFormats()
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '233.00'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -2.81%  '
$ws.Range('E5').ClearFormats()
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('E6').ClearFormats()
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4689'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -2.39%  '
$ws.Range('E7').ClearFormats()
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2707'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -4.61%  '
$ws.Range('E8').ClearFormats()
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06372'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -2.73%  '
$ws.Range('E9').ClearFormats()
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.871.50'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -1.48%  '
$ws.Range('E10').ClearFormats()
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07425'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('E11').ClearFormats()
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '16.21'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -2.75%  '
$ws.Range('E12').ClearFormats()
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.950'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -3.01%  '
$ws.Range('E13').ClearFormats()
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '84.95'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -3.57%  '
$ws.Range('E14').ClearFormats()
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6266'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -6.34%  '
$ws.Range('E15').ClearFormats()
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '30.137.06'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -1.75%  '
$ws.Range('E16').ClearFormats()
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '229.82'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('E18').ClearFormats()
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.62'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -5.24%  '
$ws.Range('E19').ClearFormats()
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000007318'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -3.95%  '
$ws.Range('E20').ClearFormats()
$ws.Range('B21').Value = 'BinanceUSD'
$ws.Range('C21').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.002'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('E21').ClearFormats()
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.943'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -6.98%  '
$ws.Range('E22').ClearFormats()
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.948'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -4.49%  '
$ws.Range('E23').ClearFormats()
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.243'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.36%  '
$ws.Range('E24').ClearFormats()
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '166.19'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -2.06%  '
$ws.Range('E25').ClearFormats()
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '17.77'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -4.76%  '
$ws.Range('E26').ClearFormats()
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.866'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -4.39%  '
$ws.Range('E27').ClearFormats()
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.1036'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +3.04%  '
$ws.Range('E28').ClearFormats()
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.389'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.85%  '
$ws.Range('E29').ClearFormats()
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.102'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -6.11%  '
$ws.Range('E30').ClearFormats()
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.877'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -3.67%  '
$ws.Range('E31').ClearFormats()
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.04882'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -3.59%  '
$ws.Range('E32').ClearFormats()
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.157'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -5.52%  '
$ws.Range('E33').ClearFormats()
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7123'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -5.43%  '
$ws.Range('E34').ClearFormats()
$ws.Range('B35').Value = 'Frax'
$ws.Range('C35').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.000'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('E35').ClearFormats()
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.703'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.40%  '
$ws.Range('E36').ClearFormats()
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.01859'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -1.30%  '
$ws.Range('E37').ClearFormats()
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.637'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('E38').ClearFormats()
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.9072'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -1.78%  '
$ws.Range('E39').ClearFormats()
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.933'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -7.28%  '
$ws.Range('E40').ClearFormats()
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '105.28'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -1.67%  '
$ws.Range('E41').ClearFormats()
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.9990'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('E42').ClearFormats()
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.543'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -5.15%  '
$ws.Range('E43').ClearFormats()
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.4063'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -5.38%  '
$ws.Range('E44').ClearFormats()
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '7.025'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -5.11%  '
$ws.Range('E45').ClearFormats()
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '60.33'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -6.57%  '
$ws.Range('E46').ClearFormats()
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.1186'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -6.86%  '
$ws.Range('E47').ClearFormats()
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.631'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -3.96%  '
$ws.Range('E48').ClearFormats()
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '33.02'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -2.84%  '
$ws.Range('E49').ClearFormats()
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.386'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -6.73%  '
$ws.Range('E50').ClearFormats()
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05577'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -1.52%  '
$ws.Range('E51').ClearFormats()
